# Update "Daily Doodh" worksheet for Nov 25, 2021 (row 26).
# Fills in the day's quantities for every supplier/buyer column (H:AY),
# which were previously blank. Downstream SUM/formula cells (BA26:BD26,
# the monthly totals in row 32-37, etc.) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dayValues = @{
    "H26"  = 0.5;
    "I26"  = 3;
    "J26"  = 1.5;
    "K26"  = 2;
    "L26"  = 0;
    "M26"  = 2;
    "N26"  = 3.5;
    "O26"  = 2.5;
    "P26"  = 1;
    "Q26"  = 5;
    "R26"  = 2;
    "S26"  = 3;
    "T26"  = 3;
    "U26"  = 0;
    "V26"  = 1.5;
    "W26"  = 1;
    "X26"  = 1;
    "Y26"  = 1;
    "Z26"  = 1;
    "AA26" = 1.5;
    "AB26" = 5;
    "AC26" = 2.5;
    "AD26" = 0;
    "AE26" = 0.5;
    "AF26" = 2;
    "AG26" = 2;
    "AH26" = 1;
    "AI26" = 0;
    "AJ26" = 1.5;
    "AK26" = 2;
    "AL26" = 3;
    "AM26" = 3.5;
    "AN26" = 1.5;
    "AO26" = 0;
    "AP26" = 1;
    "AQ26" = 1.5;
    "AR26" = 1.5;
    "AS26" = 0;
    "AT26" = 0;
    "AU26" = 2.5;
    "AV26" = 0.5;
    "AW26" = 2.5;
    "AX26" = 0;
    "AY26" = 2
}

foreach ($cellRef in $dayValues.Keys) {
    $ws.Range($cellRef).Value = $dayValues[$cellRef]
}

# Update the view state to match: scrolled right to column AR, with BB26
# as the active/selected cell.
[void]$ws.Range("BB26").Select()
$excel.ActiveWindow.ScrollColumn = 44
$excel.ActiveWindow.ScrollRow = 1

Write-Output "Updated Nov 25, 2021 (row 26) daily doodh entries."
